# EverNoteApi_map.xlsx - "added Tag Controller Integration tests"
#
# Reworks the /person/{id}/pad and /person/{id}/pad/.../note sections to
# use {padId} instead of {name}/{padName}, tightens a few Result
# descriptions, and adds a brand new "tag" sub-resource nested under
# /person/{id}/pad/{padId}/note/{noteId}, plus restates the existing
# top-level /person/{id}/tag section right below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# the spacer row between the note and tag blocks no longer carries the
# old section's row format
$ws.Rows.Item(11).ClearFormats()

# ---- /person/{id}/pad section (rows 6-7) ----------------------------
$ws.Range("C6").Value2 = "Get all the notepads for a person"

$ws.Range("A7").Value2 = "/person/{id}/pad/{padId}"
$ws.Range("C7").Value2 = "Get, create, update or delete single pad"

# ---- /person/{id}/pad/{padId}/note section (rows 9-10) ---------------
$ws.Range("A9").Value2 = "/person/{id}/pad/{padId}/note"
$ws.Range("C9").Value2 = "Get all the notes for a pad"

$ws.Range("A10").Value2 = "/person/{id}/pad/{padId}/note/{noteId}"
$ws.Range("C10").Value2 = "Get, create, update or delete single note"

# ---- NEW: /person/{id}/pad/{padId}/note/{noteId}/tag section ---------
# (reuses rows 12-13, which previously held the old top-level tag rows)
$ws.Range("A12").Value2 = "/person/{id}/pad/{padId}/note/{noteId}/tag"
$ws.Range("B12").Value2 = "GET"
$ws.Range("C12").Value2 = "Get all the tags for a note"

$ws.Range("A13").Value2 = "/person/{id}/pad/{padId}/note/{noteId}/tag/{tagName}"
$ws.Range("B13").Value2 = "POST,DELETE"
$ws.Range("C13").Value2 = "Create or delete single tag for a note"

# ---- restated top-level /person/{id}/tag section (new rows 15-16) ----
$ws.Range("A15").Value2 = "/person/{id}/tag"
$ws.Range("B15").Value2 = "GET"
$ws.Range("C15").Value2 = "Get all the tags for a person"

$ws.Range("A16").Value2 = "/person/{id}/tag/{tagName}"
$ws.Range("B16").Value2 = "GET"
$ws.Range("C16").Value2 = "Get all the tags with name {tagName} for a person"

# ---- selection left where the author last clicked after editing -----
$ws.Range("A18").Select()
